$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.201.37'
$ws.Range("E2").Value = '  -0.50%  '

# Row 3
$ws.Range("D3").Value = '3.425.96'
$ws.Range("E3").Value = '  -0.92%  '

# Row 4
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '408.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.49%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.75%  '

# Row 7
$ws.Range("E7").Value = '  +0.28%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.15%  '

# Row 9
$ws.Range("E9").Value = '  -1.60%  '

# Row 10
$ws.Range("E10").Value = '  -3.12%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.28'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.17%  '

# Row 12
$ws.Range("E12").Value = '  -1.57%  '

# Row 13
$ws.Range("D13").Value = '3.969.57'
$ws.Range("E13").Value = '  -1.01%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.15%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.98'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.46%  '

# Row 16
$ws.Range("D16").Value = '3.421.60'
$ws.Range("E16").Value = '  -1.01%  '

# Row 17
$ws.Range("D17").Value = '62.203.93'
$ws.Range("E17").Value = '  -0.42%  '

# Row 18
$ws.Range("E18").Value = '  -2.41%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.98%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000133'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.02%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '85.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '315.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.70%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.88%  '

# Row 25
$ws.Range("E25").Value = '  -2.50%  '

# Row 26
$ws.Range("E26").Value = '  +9.45%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.55%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.81%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.82%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.21%  '

# Row 31
$ws.Range("E31").Value = '  -1.85%  '

# Row 32
$ws.Range("E32").Value = '  -5.33%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '42.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.44%  '

# Row 34
$ws.Range("E34").Value = '  -4.29%  '

# Row 35
$ws.Range("E35").Value = '  +0.05%  '

# Row 36
$ws.Range("E36").Value = '  -2.00%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.12%  '

# Row 39
$ws.Range("E39").Value = '  -3.86%  '

# Row 40
$ws.Range("E40").Value = '  -0.95%  '

# Row 41
$ws.Range("E41").Value = '  +0.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.73%  '

# Row 43
$ws.Range("E43").Value = '  -0.67%  '

# Row 44
$ws.Range("E44").Value = '  +1.36%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.04%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.89'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.78%  '

# Row 47
$ws.Range("E47").Value = '  -3.15%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.38%  '

# Row 49
$ws.Range("D49").Value = '2.132.19'
$ws.Range("E49").Value = '  -5.16%  '

# Row 50
$ws.Range("E50").Value = '  -4.42%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.96%  '
